$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 27.03890566666666
$ws.Range("H2").Value = 81.116717
$ws.Range("I2").Value = 0.07096188219033728
$ws.Range("J2").Value = 0.07096188219033729
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7065936666666666
$ws.Range("N2").Value = 2.119781
$ws.Range("O2").Value = 0.005187843618793344
$ws.Range("P2").Value = 0.005187843618793344
$ws.Range("Q2").Value = 19.10551949766411
$ws.Range("R2").Value = 171.949675478977
$ws.Range("S2").Value = 0.0003681391476987063
$ws.Range("T2").Value = 0.0003681391476987064

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 27.03890566666666
$ws.Range("H3").Value = 81.116717
$ws.Range("I3").Value = 0.07096188219033728
$ws.Range("J3").Value = 0.07096188219033729
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 111.9320066666667
$ws.Range("N3").Value = 335.79602
$ws.Range("O3").Value = 0.8218100075305903
$ws.Range("P3").Value = 0.8218100075305903
$ws.Range("Q3").Value = 3026.518969340704
$ws.Range("R3").Value = 27238.67072406634
$ws.Range("S3").Value = 0.05831718493722594
$ws.Range("T3").Value = 0.05831718493722595

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 27.03890566666666
$ws.Range("H4").Value = 81.116717
$ws.Range("I4").Value = 0.07096188219033728
$ws.Range("J4").Value = 0.07096188219033729
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 23.563205
$ws.Range("N4").Value = 70.689615
$ws.Range("O4").Value = 0.1730021488506163
$ws.Range("P4").Value = 0.1730021488506163
$ws.Range("Q4").Value = 637.1232771993283
$ws.Range("R4").Value = 5734.109494793955
$ws.Range("S4").Value = 0.01227655810541263
$ws.Range("T4").Value = 0.01227655810541263

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 345.566579
$ws.Range("H5").Value = 1036.699737
$ws.Range("I5").Value = 0.9069174311350353
$ws.Range("J5").Value = 0.9069174311350354
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.7065936666666666
$ws.Range("N5").Value = 2.119781
$ws.Range("O5").Value = 0.005187843618793344
$ws.Range("P5").Value = 0.005187843618793344
$ws.Range("Q5").Value = 244.1751561330663
$ws.Range("R5").Value = 2197.576405197597
$ws.Range("S5").Value = 0.004704945807886345
$ws.Range("T5").Value = 0.004704945807886346

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 345.566579
$ws.Range("H6").Value = 1036.699737
$ws.Range("I6").Value = 0.9069174311350353
$ws.Range("J6").Value = 0.9069174311350354
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 111.9320066666667
$ws.Range("N6").Value = 335.79602
$ws.Range("O6").Value = 0.8218100075305903
$ws.Range("P6").Value = 0.8218100075305903
$ws.Range("Q6").Value = 38679.96062440519
$ws.Range("R6").Value = 348119.6456196467
$ws.Range("S6").Value = 0.745313820910707
$ws.Range("T6").Value = 0.7453138209107071

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 345.566579
$ws.Range("H7").Value = 1036.699737
$ws.Range("I7").Value = 0.9069174311350353
$ws.Range("J7").Value = 0.9069174311350354
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 23.563205
$ws.Range("N7").Value = 70.689615
$ws.Range("O7").Value = 0.1730021488506163
$ws.Range("P7").Value = 0.1730021488506163
$ws.Range("Q7").Value = 8142.656142125695
$ws.Range("R7").Value = 73283.90527913126
$ws.Range("S7").Value = 0.1568986644164419
$ws.Range("T7").Value = 0.156898664416442

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 8.428738666666668
$ws.Range("H8").Value = 25.286216
$ws.Range("I8").Value = 0.0221206866746274
$ws.Range("J8").Value = 0.02212068667462741
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7065936666666666
$ws.Range("N8").Value = 2.119781
$ws.Range("O8").Value = 0.005187843618793344
$ws.Range("P8").Value = 0.005187843618793344
$ws.Range("Q8").Value = 5.955693359855111
$ws.Range("R8").Value = 53.601240238696
$ws.Range("S8").Value = 0.0001147586632082927
$ws.Range("T8").Value = 0.0001147586632082927

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 8.428738666666668
$ws.Range("H9").Value = 25.286216
$ws.Range("I9").Value = 0.0221206866746274
$ws.Range("J9").Value = 0.02212068667462741
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 111.9320066666667
$ws.Range("N9").Value = 335.79602
$ws.Range("O9").Value = 0.8218100075305903
$ws.Range("P9").Value = 0.8218100075305903
$ws.Range("Q9").Value = 943.4456326289246
$ws.Range("R9").Value = 8491.01069366032
$ws.Range("S9").Value = 0.01817900168265737
$ws.Range("T9").Value = 0.01817900168265738

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 8.428738666666668
$ws.Range("H10").Value = 25.286216
$ws.Range("I10").Value = 0.0221206866746274
$ws.Range("J10").Value = 0.02212068667462741
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 23.563205
$ws.Range("N10").Value = 70.689615
$ws.Range("O10").Value = 0.1730021488506163
$ws.Range("P10").Value = 0.1730021488506163
$ws.Range("Q10").Value = 198.6080970940934
$ws.Range("R10").Value = 1787.47287384684
$ws.Range("S10").Value = 0.003826926328761734
$ws.Range("T10").Value = 0.003826926328761736
